$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WBS")

for ($row = 2; $row -le 49; $row++) {
    $ws.Range("H$row").Value = "Done"
    $ws.Range("J$row").Copy($ws.Range("K$row"))
    $ws.Range("L$row").Value = "✅"
    $ws.Range("M$row").Value = "✅"
    $ws.Range("N$row").Value = "✅"
    $ws.Range("O$row").Value = "✅"
    $ws.Range("P$row").Value = "✅"
}
